# Avanzamento.xlsx — avanzamento update (date shift for the in-progress batch)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Column B holds the "data" (date) for each record. Rows 132-198 were all
# scheduled for 45937 (07/10/2025); push that whole block out 10 days to
# 45947 (17/10/2025).
for ($r = 132; $r -le 198; $r++) {
    $ws.Cells.Item($r, 2).Value = 45947
}

# Reflect where the user ended up looking/working: scrolled further down the
# sheet and selected B195.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 173
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B195").Select()
